$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: add Pwd@013506 to B6, "test fail" to C6
$ws.Range("B6").Value = "Pwd@013506"
$ws.Range("C6").Value = "test fail"

# Row 7: 013506 / Pwd@013506 / test fail
$ws.Range("A7").Value = "013506"
$ws.Range("B7").Value = "Pwd@013506"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C7").Value = "test fail"

# Update selection to E13
$ws.Range("E13").Select()
